$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "CV" column (G), shifting the
# existing G:K columns (CV, HEX, Error, Error Cents, note) one to the
# right (H:L).
$ws.Columns("G").Insert()

# New column header (row 5) -- "mV"
$ws.Range("G5").Value = "mV"

# Row 6: one-off (non filled-down) formula for the new mV column.
$ws.Range("G6").Formula = "=ROUND(F6*1000, 0)"

# Rows 7:18 were filled down together as one shared formula.
$ws.Range("G7:G18").Formula = "=ROUND(F7*1000, 0)"

# New data point added below the table.
$ws.Range("M9").Value = 134940

# Scroll the view back to the top-left and leave the selection on the
# newly touched cell, matching where Excel lands after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("M10").Select()
